$wb = $excel.ActiveWorkbook

# --- Sheet: team_df ---
$wsTeam = $wb.Worksheets.Item("team_df")

# Row 6: arrive_thres count increases from 5 to 6 (out of 22)
$wsTeam.Range("S6").Value = 6
$wsTeam.Range("U6").Value = 0.2727272727272727

# Row 7: arrive_thres count increases from 4 to 5 (out of 21)
$wsTeam.Range("S7").Value = 5
$wsTeam.Range("U7").Value = 0.2380952380952381

# --- Sheet: team_df_day ---
$wsDay = $wb.Worksheets.Item("team_df_day")

# Row 3: arrive_thres count increases from 8 to 9 (out of 48)
$wsDay.Range("F3").Value = 9
$wsDay.Range("H3").Value = 0.1875

# Row 7: arrive_thres count increases from 12 to 13 (out of 63)
$wsDay.Range("F7").Value = 13
$wsDay.Range("H7").Value = 0.2063492063492063

# --- Sheet: productivity_tl ---
$wsTl = $wb.Worksheets.Item("productivity_tl")

$wsTl.Range("D3").Value = 0.1875
$wsTl.Range("D7").Value = 0.2063492063492063

# --- Sheet: productivity_team_function ---
$wsFunc = $wb.Worksheets.Item("productivity_team_function")

$wsFunc.Range("D3").Value = 0.1875
$wsFunc.Range("D7").Value = 0.2063492063492063
